$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Supplier role) test account values, replacing Sean's test
# data with the supplier test data.
$ws.Range("B2").Value = "supplierTest"
$ws.Range("D2").Value = "SUPPLIER TEST"

# Move the selection to reflect the final cursor position after editing.
$ws.Range("E11").Select()
